# Cleans up vaccine/brand name labels across all worksheets:
#  - strips footnote markers like " [1]", " [2]", etc.
#  - collapses embedded line breaks (used to wrap text in a cell) into a
#    single space, joining what used to be separate lines.
#
# This mirrors the source workbook edit where labels such as
# "DTaP [1]" became "DTaP " and multi-line labels such as
# "Recombivax`nHB" became "Recombivax HB" (and similarly
# "Afluria`nQuadrivalent" became "Afluria Quadrivalent", which then
# collapses into the already-existing "Afluria Quadrivalent" entry).

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    $cols = $used.Columns.Count

    for ($r = 1; $r -le $rows; $r++) {
        for ($c = 1; $c -le $cols; $c++) {
            $cell = $used.Cells.Item($r, $c)
            $val = $cell.Value2

            if ($val -is [string]) {
                $newVal = $val -replace '\[\d+\]', ''
                $newVal = $newVal -replace "`r`n", ' '
                $newVal = $newVal -replace "`n", ' '

                if ($newVal -ne $val) {
                    $cell.Value = $newVal
                }
            }
        }
    }
}
